$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -7.999
$ws.Range("D6").Value = -8.071999999999999
$ws.Range("D7").Value = -8.029
$ws.Range("D8").Value = -8.161
$ws.Range("D16").Value = -8.451000000000001
$ws.Range("D20").Value = -8.058000000000002
$ws.Range("D21").Value = -8.337
$ws.Range("D28").Value = -7.834000000000001
$ws.Range("D29").Value = -7.56
$ws.Range("D30").Value = -7.211000000000001
$ws.Range("D32").Value = -7.276999999999999
$ws.Range("D40").Value = -8.077999999999999
$ws.Range("D46").Value = -8.022
$ws.Range("D51").Value = -8.241000000000001
$ws.Range("D52").Value = -7.904000000000001
$ws.Range("D57").Value = -7.921000000000001
$ws.Range("D59").Value = -8.058
$ws.Range("D62").Value = -7.840999999999999
$ws.Range("D66").Value = -7.102000000000001
$ws.Range("D73").Value = -8.038999999999998
$ws.Range("D74").Value = -8.154999999999998
$ws.Range("D77").Value = -7.858
$ws.Range("D92").Value = -6.738
$ws.Range("D100").Value = -7.937
